$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 120
$ws.Cells.Item(120, 4).Value = 44466
$ws.Cells.Item(120, 9).Value = "Segunda"
$ws.Cells.Item(120, 10).Value = 500

# Row 121
$ws.Cells.Item(121, 4).Value = 44389

# Row 122
$ws.Cells.Item(122, 4).Value = 44265
$ws.Cells.Item(122, 10).Value = 250

# Row 123
$ws.Cells.Item(123, 4).Value = 44343
$ws.Cells.Item(123, 10).Value = 500
$ws.Cells.Item(123, 11).Value = 1000
$ws.Cells.Item(123, 12).Value = 1000
$ws.Cells.Item(123, 13).Value = 1000
$ws.Cells.Item(123, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(123, 15).Value = "Región del Maule"
$ws.Cells.Item(123, 16).Value = 200
$ws.Cells.Item(123, 17).Value = 5

# Row 124
$ws.Cells.Item(124, 10).Value = 110
$ws.Cells.Item(124, 11).Value = 8000
$ws.Cells.Item(124, 12).Value = 8000
$ws.Cells.Item(124, 13).Value = 8000
$ws.Cells.Item(124, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(124, 15).Value = "Región Metropolitana"
$ws.Cells.Item(124, 16).Value = 533
$ws.Cells.Item(124, 17).Value = 15

# Row 125
$ws.Cells.Item(125, 4).Value = 44201
$ws.Cells.Item(125, 10).Value = 1200
$ws.Cells.Item(125, 11).Value = 900
$ws.Cells.Item(125, 12).Value = 900
$ws.Cells.Item(125, 13).Value = 900
$ws.Cells.Item(125, 16).Value = 180

# Row 126
$ws.Cells.Item(126, 4).Value = 44280
$ws.Cells.Item(126, 11).Value = 850
$ws.Cells.Item(126, 13).Value = 925
$ws.Cells.Item(126, 16).Value = 185

# Row 127
$ws.Cells.Item(127, 4).Value = 44270
$ws.Cells.Item(127, 10).Value = 500
$ws.Cells.Item(127, 11).Value = 800
$ws.Cells.Item(127, 13).Value = 900
$ws.Cells.Item(127, 16).Value = 180

# Row 128
$ws.Cells.Item(128, 4).Value = 44260
$ws.Cells.Item(128, 10).Value = 1200
$ws.Cells.Item(128, 11).Value = 1000
$ws.Cells.Item(128, 13).Value = 1000
$ws.Cells.Item(128, 16).Value = 200

# Row 129
$ws.Cells.Item(129, 4).Value = 44267
$ws.Cells.Item(129, 10).Value = 1000
$ws.Cells.Item(129, 11).Value = 800
$ws.Cells.Item(129, 13).Value = 900
$ws.Cells.Item(129, 16).Value = 180

# Row 130
$ws.Cells.Item(130, 4).Value = 44312
$ws.Cells.Item(130, 10).Value = 250
$ws.Cells.Item(130, 11).Value = 1000
$ws.Cells.Item(130, 12).Value = 1000
$ws.Cells.Item(130, 13).Value = 1000
$ws.Cells.Item(130, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(130, 15).Value = "Región del Maule"
$ws.Cells.Item(130, 16).Value = 200
$ws.Cells.Item(130, 17).Value = 5

# Row 131
$ws.Cells.Item(131, 10).Value = 120
$ws.Cells.Item(131, 11).Value = 10000
$ws.Cells.Item(131, 12).Value = 10000
$ws.Cells.Item(131, 13).Value = 10000
$ws.Cells.Item(131, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(131, 15).Value = "Región Metropolitana"
$ws.Cells.Item(131, 16).Value = 667
$ws.Cells.Item(131, 17).Value = 15

# Row 132
$ws.Cells.Item(132, 4).Value = 44390
$ws.Cells.Item(132, 10).Value = 1200
$ws.Cells.Item(132, 11).Value = 1000
$ws.Cells.Item(132, 12).Value = 1000
$ws.Cells.Item(132, 13).Value = 1000
$ws.Cells.Item(132, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(132, 15).Value = "Región del Maule"
$ws.Cells.Item(132, 16).Value = 200
$ws.Cells.Item(132, 17).Value = 5

# Row 133
$ws.Cells.Item(133, 10).Value = 120
$ws.Cells.Item(133, 11).Value = 10000
$ws.Cells.Item(133, 12).Value = 10000
$ws.Cells.Item(133, 13).Value = 10000
$ws.Cells.Item(133, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(133, 15).Value = "Región Metropolitana"
$ws.Cells.Item(133, 16).Value = 667
$ws.Cells.Item(133, 17).Value = 15

# Row 134
$ws.Cells.Item(134, 4).Value = 44386
$ws.Cells.Item(134, 10).Value = 1200

# Row 135
$ws.Cells.Item(135, 4).Value = 44264
$ws.Cells.Item(135, 10).Value = 1000
$ws.Cells.Item(135, 11).Value = 1000
$ws.Cells.Item(135, 12).Value = 1000
$ws.Cells.Item(135, 13).Value = 1000
$ws.Cells.Item(135, 16).Value = 200

# Row 136
$ws.Cells.Item(136, 4).Value = 44463
$ws.Cells.Item(136, 10).Value = 1200
$ws.Cells.Item(136, 11).Value = 1200
$ws.Cells.Item(136, 12).Value = 1200
$ws.Cells.Item(136, 13).Value = 1200
$ws.Cells.Item(136, 16).Value = 240

# Row 137
$ws.Cells.Item(137, 4).Value = 44301
$ws.Cells.Item(137, 10).Value = 500

# Row 138
$ws.Cells.Item(138, 4).Value = 44243
$ws.Cells.Item(138, 10).Value = 1200

# Row 139
$ws.Cells.Item(139, 4).Value = 44252
$ws.Cells.Item(139, 10).Value = 750
$ws.Cells.Item(139, 11).Value = 1000
$ws.Cells.Item(139, 12).Value = 1000
$ws.Cells.Item(139, 13).Value = 1000
$ws.Cells.Item(139, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(139, 15).Value = "Región del Maule"
$ws.Cells.Item(139, 16).Value = 200
$ws.Cells.Item(139, 17).Value = 5

# Row 140
$ws.Cells.Item(140, 10).Value = 110
$ws.Cells.Item(140, 11).Value = 7000
$ws.Cells.Item(140, 12).Value = 7000
$ws.Cells.Item(140, 13).Value = 7000
$ws.Cells.Item(140, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(140, 15).Value = "Región Metropolitana"
$ws.Cells.Item(140, 16).Value = 467
$ws.Cells.Item(140, 17).Value = 15

# Row 141
$ws.Cells.Item(141, 4).Value = 44166
$ws.Cells.Item(141, 10).Value = 1200
$ws.Cells.Item(141, 12).Value = 1000
$ws.Cells.Item(141, 13).Value = 925
$ws.Cells.Item(141, 16).Value = 185

# Row 142
$ws.Cells.Item(142, 4).Value = 44168
$ws.Cells.Item(142, 10).Value = 500
$ws.Cells.Item(142, 11).Value = 850
$ws.Cells.Item(142, 12).Value = 850
$ws.Cells.Item(142, 13).Value = 850
$ws.Cells.Item(142, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(142, 15).Value = "Región del Maule"
$ws.Cells.Item(142, 16).Value = 170
$ws.Cells.Item(142, 17).Value = 5

# Row 143
$ws.Cells.Item(143, 10).Value = 120
$ws.Cells.Item(143, 11).Value = 9000
$ws.Cells.Item(143, 12).Value = 9000
$ws.Cells.Item(143, 13).Value = 9000
$ws.Cells.Item(143, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(143, 15).Value = "Región Metropolitana"
$ws.Cells.Item(143, 16).Value = 600
$ws.Cells.Item(143, 17).Value = 15

# Row 144
$ws.Cells.Item(144, 4).Value = 44369
$ws.Cells.Item(144, 10).Value = 1400
$ws.Cells.Item(144, 11).Value = 1000
$ws.Cells.Item(144, 12).Value = 1000
$ws.Cells.Item(144, 13).Value = 1000
$ws.Cells.Item(144, 16).Value = 200

# Row 145
$ws.Cells.Item(145, 4).Value = 44221
$ws.Cells.Item(145, 11).Value = 900
$ws.Cells.Item(145, 12).Value = 900
$ws.Cells.Item(145, 13).Value = 900
$ws.Cells.Item(145, 16).Value = 180

# Row 146
$ws.Cells.Item(146, 4).Value = 44371
$ws.Cells.Item(146, 10).Value = 500
$ws.Cells.Item(146, 11).Value = 1000
$ws.Cells.Item(146, 12).Value = 1000
$ws.Cells.Item(146, 13).Value = 1000
$ws.Cells.Item(146, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(146, 15).Value = "Región del Maule"
$ws.Cells.Item(146, 16).Value = 200
$ws.Cells.Item(146, 17).Value = 5

# Row 147
$ws.Cells.Item(147, 10).Value = 110
$ws.Cells.Item(147, 11).Value = 8000
$ws.Cells.Item(147, 12).Value = 8000
$ws.Cells.Item(147, 13).Value = 8000
$ws.Cells.Item(147, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(147, 15).Value = "Región Metropolitana"
$ws.Cells.Item(147, 16).Value = 533
$ws.Cells.Item(147, 17).Value = 15

# Row 148
$ws.Cells.Item(148, 4).Value = 44316

# Row 149
$ws.Cells.Item(149, 4).Value = 44397
$ws.Cells.Item(149, 10).Value = 1200

# Row 150
$ws.Cells.Item(150, 4).Value = 44363
$ws.Cells.Item(150, 10).Value = 100
$ws.Cells.Item(150, 11).Value = 1000
$ws.Cells.Item(150, 12).Value = 1000
$ws.Cells.Item(150, 13).Value = 1000
$ws.Cells.Item(150, 16).Value = 200

# Row 151
$ws.Cells.Item(151, 4).Value = 44277
$ws.Cells.Item(151, 10).Value = 500

# Row 152
$ws.Cells.Item(152, 4).Value = 44291
$ws.Cells.Item(152, 10).Value = 250
$ws.Cells.Item(152, 12).Value = 850
$ws.Cells.Item(152, 13).Value = 850
$ws.Cells.Item(152, 16).Value = 170

# Row 153
$ws.Cells.Item(153, 4).Value = 44273
$ws.Cells.Item(153, 11).Value = 850
$ws.Cells.Item(153, 13).Value = 925
$ws.Cells.Item(153, 16).Value = 185

# Row 154
$ws.Cells.Item(154, 4).Value = 44438
$ws.Cells.Item(154, 11).Value = 1000
$ws.Cells.Item(154, 12).Value = 1000
$ws.Cells.Item(154, 13).Value = 1000
$ws.Cells.Item(154, 16).Value = 200

# Row 155
$ws.Cells.Item(155, 4).Value = 44209
$ws.Cells.Item(155, 10).Value = 500
$ws.Cells.Item(155, 11).Value = 900
$ws.Cells.Item(155, 12).Value = 900
$ws.Cells.Item(155, 13).Value = 900
$ws.Cells.Item(155, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(155, 15).Value = "Región del Maule"
$ws.Cells.Item(155, 16).Value = 180
$ws.Cells.Item(155, 17).Value = 5

# Row 156
$ws.Cells.Item(156, 4).Value = 44160
$ws.Cells.Item(156, 10).Value = 20
$ws.Cells.Item(156, 11).Value = 8000
$ws.Cells.Item(156, 12).Value = 8000
$ws.Cells.Item(156, 13).Value = 8000
$ws.Cells.Item(156, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(156, 15).Value = "Región Metropolitana"
$ws.Cells.Item(156, 16).Value = 533
$ws.Cells.Item(156, 17).Value = 15

# Row 157
$ws.Cells.Item(157, 4).Value = 44351
$ws.Cells.Item(157, 10).Value = 1200
$ws.Cells.Item(157, 11).Value = 1000
$ws.Cells.Item(157, 12).Value = 1000
$ws.Cells.Item(157, 13).Value = 1000
$ws.Cells.Item(157, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(157, 15).Value = "Región del Maule"
$ws.Cells.Item(157, 16).Value = 200
$ws.Cells.Item(157, 17).Value = 5

# Row 158
$ws.Cells.Item(158, 10).Value = 120
$ws.Cells.Item(158, 11).Value = 9000
$ws.Cells.Item(158, 12).Value = 9000
$ws.Cells.Item(158, 13).Value = 9000
$ws.Cells.Item(158, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(158, 15).Value = "Región Metropolitana"
$ws.Cells.Item(158, 16).Value = 600
$ws.Cells.Item(158, 17).Value = 15

# Row 159
$ws.Cells.Item(159, 4).Value = 44365

# Row 160
$ws.Cells.Item(160, 4).Value = 44306
$ws.Cells.Item(160, 10).Value = 1200
$ws.Cells.Item(160, 11).Value = 1000
$ws.Cells.Item(160, 12).Value = 1000
$ws.Cells.Item(160, 13).Value = 1000
$ws.Cells.Item(160, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(160, 15).Value = "Región del Maule"
$ws.Cells.Item(160, 16).Value = 200
$ws.Cells.Item(160, 17).Value = 5

# Row 161
$ws.Cells.Item(161, 10).Value = 150
$ws.Cells.Item(161, 11).Value = 10000
$ws.Cells.Item(161, 12).Value = 10000
$ws.Cells.Item(161, 13).Value = 10000
$ws.Cells.Item(161, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(161, 15).Value = "Región Metropolitana"
$ws.Cells.Item(161, 16).Value = 667
$ws.Cells.Item(161, 17).Value = 15

# Row 162
$ws.Cells.Item(162, 4).Value = 44215
$ws.Cells.Item(162, 10).Value = 1200
$ws.Cells.Item(162, 11).Value = 900
$ws.Cells.Item(162, 12).Value = 1000
$ws.Cells.Item(162, 13).Value = 950
$ws.Cells.Item(162, 16).Value = 190

# Row 163
$ws.Cells.Item(163, 4).Value = 44175
$ws.Cells.Item(163, 10).Value = 500
$ws.Cells.Item(163, 11).Value = 850
$ws.Cells.Item(163, 12).Value = 850
$ws.Cells.Item(163, 13).Value = 850
$ws.Cells.Item(163, 16).Value = 170

# Row 164
$ws.Cells.Item(164, 4).Value = 44357
$ws.Cells.Item(164, 10).Value = 400
$ws.Cells.Item(164, 11).Value = 1000
$ws.Cells.Item(164, 12).Value = 1000
$ws.Cells.Item(164, 13).Value = 1000
$ws.Cells.Item(164, 16).Value = 200

# Row 165
$ws.Cells.Item(165, 4).Value = 44203
$ws.Cells.Item(165, 10).Value = 500
$ws.Cells.Item(165, 11).Value = 850
$ws.Cells.Item(165, 12).Value = 900
$ws.Cells.Item(165, 13).Value = 875
$ws.Cells.Item(165, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(165, 15).Value = "Región del Maule"
$ws.Cells.Item(165, 16).Value = 175
$ws.Cells.Item(165, 17).Value = 5

# Row 166
$ws.Cells.Item(166, 10).Value = 120
$ws.Cells.Item(166, 11).Value = 8000
$ws.Cells.Item(166, 12).Value = 8000
$ws.Cells.Item(166, 13).Value = 8000
$ws.Cells.Item(166, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(166, 15).Value = "Región Metropolitana"
$ws.Cells.Item(166, 16).Value = 533
$ws.Cells.Item(166, 17).Value = 15

# Row 167
$ws.Cells.Item(167, 4).Value = 44162
$ws.Cells.Item(167, 10).Value = 1200

# Row 168
$ws.Cells.Item(168, 4).Value = 44410
$ws.Cells.Item(168, 10).Value = 500

# Row 169
$ws.Cells.Item(169, 4).Value = 44411

# Row 170
$ws.Cells.Item(170, 4).Value = 44257
$ws.Cells.Item(170, 10).Value = 1200

# Row 171
$ws.Cells.Item(171, 4).Value = 44244
$ws.Cells.Item(171, 10).Value = 250
$ws.Cells.Item(171, 11).Value = 1000
$ws.Cells.Item(171, 12).Value = 1000
$ws.Cells.Item(171, 13).Value = 1000
$ws.Cells.Item(171, 16).Value = 200

# Row 172
$ws.Cells.Item(172, 4).Value = 44176
$ws.Cells.Item(172, 11).Value = 850
$ws.Cells.Item(172, 12).Value = 850
$ws.Cells.Item(172, 13).Value = 850
$ws.Cells.Item(172, 16).Value = 170

# Row 173
$ws.Cells.Item(173, 4).Value = 44239
$ws.Cells.Item(173, 10).Value = 1200
$ws.Cells.Item(173, 11).Value = 1000
$ws.Cells.Item(173, 12).Value = 1000
$ws.Cells.Item(173, 13).Value = 1000
$ws.Cells.Item(173, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(173, 15).Value = "Región del Maule"
$ws.Cells.Item(173, 16).Value = 200
$ws.Cells.Item(173, 17).Value = 5

# Row 174
$ws.Cells.Item(174, 10).Value = 110
$ws.Cells.Item(174, 11).Value = 9000
$ws.Cells.Item(174, 12).Value = 9000
$ws.Cells.Item(174, 13).Value = 9000
$ws.Cells.Item(174, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(174, 15).Value = "Región Metropolitana"
$ws.Cells.Item(174, 16).Value = 600
$ws.Cells.Item(174, 17).Value = 15

# Row 175
$ws.Cells.Item(175, 4).Value = 44376
$ws.Cells.Item(175, 10).Value = 1200
$ws.Cells.Item(175, 11).Value = 1000
$ws.Cells.Item(175, 12).Value = 1000
$ws.Cells.Item(175, 13).Value = 1000
$ws.Cells.Item(175, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(175, 15).Value = "Región del Maule"
$ws.Cells.Item(175, 16).Value = 200
$ws.Cells.Item(175, 17).Value = 5

# Row 176
$ws.Cells.Item(176, 10).Value = 120
$ws.Cells.Item(176, 11).Value = 7000
$ws.Cells.Item(176, 12).Value = 8000
$ws.Cells.Item(176, 13).Value = 7500
$ws.Cells.Item(176, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(176, 15).Value = "Región Metropolitana"
$ws.Cells.Item(176, 16).Value = 500
$ws.Cells.Item(176, 17).Value = 15

# Row 177
$ws.Cells.Item(177, 4).Value = 44292
$ws.Cells.Item(177, 10).Value = 1000
$ws.Cells.Item(177, 11).Value = 850
$ws.Cells.Item(177, 12).Value = 1000
$ws.Cells.Item(177, 13).Value = 925
$ws.Cells.Item(177, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(177, 15).Value = "Región del Maule"
$ws.Cells.Item(177, 16).Value = 185
$ws.Cells.Item(177, 17).Value = 5

# Row 178
$ws.Cells.Item(178, 10).Value = 120
$ws.Cells.Item(178, 11).Value = 8000
$ws.Cells.Item(178, 12).Value = 8000
$ws.Cells.Item(178, 13).Value = 8000
$ws.Cells.Item(178, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(178, 15).Value = "Región Metropolitana"
$ws.Cells.Item(178, 16).Value = 533
$ws.Cells.Item(178, 17).Value = 15

# Row 179
$ws.Cells.Item(179, 4).Value = 44358
$ws.Cells.Item(179, 10).Value = 1000
$ws.Cells.Item(179, 11).Value = 1000
$ws.Cells.Item(179, 12).Value = 1000
$ws.Cells.Item(179, 13).Value = 1000
$ws.Cells.Item(179, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(179, 15).Value = "Región del Maule"
$ws.Cells.Item(179, 16).Value = 200
$ws.Cells.Item(179, 17).Value = 5

# Row 180
$ws.Cells.Item(180, 10).Value = 120
$ws.Cells.Item(180, 11).Value = 10000
$ws.Cells.Item(180, 12).Value = 10000
$ws.Cells.Item(180, 13).Value = 10000
$ws.Cells.Item(180, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(180, 15).Value = "Región Metropolitana"
$ws.Cells.Item(180, 16).Value = 667
$ws.Cells.Item(180, 17).Value = 15

# Row 181
$ws.Cells.Item(181, 4).Value = 44211
$ws.Cells.Item(181, 10).Value = 1000
$ws.Cells.Item(181, 11).Value = 900
$ws.Cells.Item(181, 13).Value = 950
$ws.Cells.Item(181, 16).Value = 190

# Row 182
$ws.Cells.Item(182, 4).Value = 44425
$ws.Cells.Item(182, 10).Value = 1200
$ws.Cells.Item(182, 11).Value = 1000
$ws.Cells.Item(182, 12).Value = 1000
$ws.Cells.Item(182, 13).Value = 1000
$ws.Cells.Item(182, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(182, 15).Value = "Región del Maule"
$ws.Cells.Item(182, 16).Value = 200
$ws.Cells.Item(182, 17).Value = 5

# Row 183
$ws.Cells.Item(183, 10).Value = 120
$ws.Cells.Item(183, 11).Value = 9000
$ws.Cells.Item(183, 12).Value = 9000
$ws.Cells.Item(183, 13).Value = 9000
$ws.Cells.Item(183, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(183, 15).Value = "Región Metropolitana"
$ws.Cells.Item(183, 16).Value = 600
$ws.Cells.Item(183, 17).Value = 15

# New row 184
$ws.Cells.Item(184, 1).Value = 4
$ws.Cells.Item(184, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(184, 3).Value = "Los Lagos"
$ws.Cells.Item(184, 4).Value = 44323
$ws.Cells.Item(184, 5).Value = 10
$ws.Cells.Item(184, 6).Value = 100114014
$ws.Cells.Item(184, 7).Value = "Betarraga"
$ws.Cells.Item(184, 8).Value = "Sin especificar"
$ws.Cells.Item(184, 9).Value = "Primera"
$ws.Cells.Item(184, 10).Value = 1000
$ws.Cells.Item(184, 11).Value = 1000
$ws.Cells.Item(184, 12).Value = 1000
$ws.Cells.Item(184, 13).Value = 1000
$ws.Cells.Item(184, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(184, 15).Value = "Región del Maule"
$ws.Cells.Item(184, 16).Value = 200
$ws.Cells.Item(184, 17).Value = 5
$ws.Cells.Item(184, 18).Value = "Hortaliza"
$ws.Cells.Item(184, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Host "done"